$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.442.35'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').Value = '3.286.08'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = "'582.86"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.70%  '
$ws.Range('D6').Value = "'181.80"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'0.587"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('D9').Value = '3.276.91'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('E10').Value = '  +2.07%  '
$ws.Range('D11').Value = "'0.575"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').Value = "'46.15"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = "'0.0000274"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.46%  '
$ws.Range('D14').Value = "'635.79"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +10.97%  '
$ws.Range('D15').Value = '3.818.98'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = "'8.40"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '67.472.95'
$ws.Range('E17').Value = '  +2.66%  '
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').Value = '3.296.84'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('D20').Value = "'17.55"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'10.86"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('D22').Value = "'0.892"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('D23').Value = "'17.68"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').Value = "'97.40"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('E26').Value = '  +1.16%  '
$ws.Range('E27').Value = '  +4.08%  '
$ws.Range('D28').Value = "'9.55"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.66%  '
$ws.Range('D29').Value = "'32.59"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.16%  '
$ws.Range('D30').Value = "'8.51"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').Value = "'6.64"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').Value = "'592.16"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.48%  '
$ws.Range('D33').Value = '3.925.88'
$ws.Range('E33').Value = '  +4.35%  '
$ws.Range('D34').Value = "'3.60"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('D35').Value = "'10.89"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('D36').Value = "'0.104"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('D37').Value = "'0.996"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = "'55.66"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = "'3.25"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.75%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = "'2.69"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.10%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = "'0.128"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.73%  '
$ws.Range('D42').Value = "'32.75"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.45%  '
$ws.Range('D43').Value = '0.0₃0684'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('D44').Value = "'3.34"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.44%  '
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('D46').Value = "'0.0412"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.47%  '
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('D50').Value = "'1.33"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.97%  '
$ws.Range('D51').Value = "'130.33"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.33%  '
